$wb = $excel.ActiveWorkbook

# --- Update "Metadata" sheet: bump the "Last Updated" timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "11 Nov 2025, 09:18 AM"

# --- Update "Industry Analysis" sheet: insert a new "indices" row at row 38 ---
$ws = $wb.Worksheets.Item("Industry Analysis")

# Insert a new blank row at position 38; this shifts existing rows 38-76 down to 39-77
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the "indices" data
$ws.Range("A38").Value = "🏭"
$ws.Range("B38").Value = "indices"
$ws.Range("C38").Value = 0.1668
$ws.Range("D38").Value = 0.4892
$ws.Range("E38").Value = -0.4086
$ws.Range("F38").Value = 8.3409
$ws.Range("G38").Value = 8.8048
$ws.Range("H38").Value = 20.0278
$ws.Range("I38").Value = 24.1192
$ws.Range("J38").Value = 14.4364
$ws.Range("K38").Value = 11.4954

# The insert pushed the former last row (76, "pharmaceuticals - indian - formulations")
# down to row 77, dropping it off the bottom of the table; remove that now-duplicate row
$ws.Rows.Item(77).Delete()
